# Weekly update: add two new "Membrillo" price rows for Terminal La Palmera
# de La Serena (new reporting date 2022-04-14 / serial 44665), inserted right
# after the existing header/first data row (row 11), pushing all subsequent
# rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 12 - everything that was row 12.. shifts to row 14..
$ws.Rows.Item(12).Resize(2).Insert()

# New row 12: "Primera" quality
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "Terminal La Palmera de La Serena"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44665
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100104
$ws.Range("H12").Value = "Frutos de pepita"
$ws.Range("I12").Value = 100104003
$ws.Range("J12").Value = "Membrillo"
$ws.Range("K12").Value = "Champion"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 10
$ws.Range("N12").Value = 300000
$ws.Range("O12").Value = 310000
$ws.Range("P12").Value = 305000
$ws.Range("Q12").Value = "$/bins (450 kilos)"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 678
$ws.Range("T12").Value = 450

# New row 13: "Segunda" quality
$ws.Range("A13").Value = 8
$ws.Range("B13").Value = "Terminal La Palmera de La Serena"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44665
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100104
$ws.Range("H13").Value = "Frutos de pepita"
$ws.Range("I13").Value = 100104003
$ws.Range("J13").Value = "Membrillo"
$ws.Range("K13").Value = "Champion"
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 20
$ws.Range("N13").Value = 270000
$ws.Range("O13").Value = 280000
$ws.Range("P13").Value = 275000
$ws.Range("Q13").Value = "$/bins (450 kilos)"
$ws.Range("R13").Value = "Región de O'Higgins"
$ws.Range("S13").Value = 611
$ws.Range("T13").Value = 450
